{"js": "// 1) Update the \"Date\" paragraph text.\nconst dateResults = context.document.body.search(\n  \"October  11, 2021 (10:26:58 PM)\",\n  { matchCase: true }\n);\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\n    \"October  12, 2021 (01:47:40 AM)\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Expand the \"first value ... in both arrays\" question into the more\n//    detailed wording that references the chars1 / chars2 arrays by name,\n//    with chars1 / chars2 styled using the \"NormalTok\" character style\n//    (matching the inline-code styling used elsewhere in the document).\nconst qResults = context.document.body.search(\n  \"What is the first value that occurs in both arrays, searching from left to right? If none is found, display\",\n  { matchCase: true }\n);\nqResults.load(\"text\");\nawait context.sync();\n\nif (qResults.items.length > 0) {\n  const target = qResults.items[0];\n\n  // Collapse the old run's text to an insertion point.\n  target.insertText(\"\", \"Replace\");\n  await context.sync();\n\n  // Each entry is [text, characterStyleNameOrNull]. Insert them in reverse\n  // order using \"Before\" on the same collapsed point so the final reading\n  // order is left-to-right and each piece becomes its own run (mirroring\n  // the run-per-token structure already used throughout this document).\n  const pieces = [\n    [\"What is the first value of the array\", null],\n    [\" \", null],\n    [\"chars1\", \"NormalTok\"],\n    [\" \", null],\n    [\"that also occurs in the second array\", null],\n    [\" \", null],\n    [\"chars2\", \"NormalTok\"],\n    [\", searching from left to right? If none is found, display\", null],\n  ];\n\n  for (let i = pieces.length - 1; i >= 0; i--) {\n    const [text, style] = pieces[i];\n    const inserted = target.insertText(text, \"Before\");\n    if (style) {\n      inserted.style = style;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the \"Date\" paragraph text.\n$dateFind = $d.Content.Find\n$dateFind.Text = \"October  11, 2021 (10:26:58 PM)\"\n$dateFind.Execute() | Out-Null\nif ($dateFind.Found) {\n    $dateFind.Parent.Text = \"October  12, 2021 (01:47:40 AM)\"\n}\n\n# 2) Expand the \"first value ... in both arrays\" question into the more\n#    detailed wording that references the chars1 / chars2 arrays by name,\n#    with chars1 / chars2 styled using the \"NormalTok\" character style\n#    (matching the inline-code styling used elsewhere in the document).\n$qFind = $d.Content.Find\n$qFind.Text = \"What is the first value that occurs in both arrays, searching from left to right? If none is found, display\"\n$qFind.Execute() | Out-Null\nif ($qFind.Found) {\n    $rng = $qFind.Parent\n    $paraStart = $rng.Start\n\n    # Collapse the old run's text to an insertion point.\n    $rng.Text = \"\"\n\n    # Insert the pieces in reverse order using InsertBefore on the same\n    # collapsed point so the final reading order is left-to-right and each\n    # piece becomes its own run (mirroring the run-per-token structure\n    # already used throughout this document).\n    $rng.InsertBefore(\", searching from left to right? If none is found, display\")\n    $rng.InsertBefore(\"chars2\")\n    $rng.InsertBefore(\" \")\n    $rng.InsertBefore(\"that also occurs in the second array\")\n    $rng.InsertBefore(\" \")\n    $rng.InsertBefore(\"chars1\")\n    $rng.InsertBefore(\" \")\n    $rng.InsertBefore(\"What is the first value of the array\")\n\n    $newText = \"What is the first value of the array chars1 that also occurs in the second array chars2, searching from left to right? If none is found, display\"\n    $paraEnd = $paraStart + $newText.Length\n\n    $c1Range = $d.Range($paraStart, $paraEnd)\n    $c1Find = $c1Range.Find\n    $c1Find.Text = \"chars1\"\n    $c1Find.Execute() | Out-Null\n    if ($c1Find.Found) {\n        $c1Range.Style = \"NormalTok\"\n    }\n\n    $c2Range = $d.Range($paraStart, $paraEnd)\n    $c2Find = $c2Range.Find\n    $c2Find.Text = \"chars2\"\n    $c2Find.Execute() | Out-Null\n    if ($c2Find.Found) {\n        $c2Range.Style = \"NormalTok\"\n    }\n}\n"}
